# Add a new "Save" column (H) to the sheet:
#  - H1 header "Save", styled the same as the other header cells (bold/bordered)
#  - H2 = 1, H3 = 0 (plain numeric values, no special style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the neighboring header cell (G1) onto the new
# header cell (H1) so it matches the other header cells' style (bold,
# centered/top aligned, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header and data values.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
